$d = $word.ActiveDocument

# Locate the "Donnees centrees reduites" Heading1 paragraph; the new Variance / Correlation
# sections must be inserted immediately before it (right after the Covariance section).
$searchRange = $d.Content.Duplicate
$found = $searchRange.Find.Execute("Données centrées réduites", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

$targetPara = $searchRange.Paragraphs(1)
$insertionPoint = $targetPara.Range.Duplicate
$insertionPoint.Collapse(1)

$newSectionsXml = '<w:p><w:pPr><w:pStyle w:val="Heading2"/></w:pPr><w:bookmarkStart w:id="27" w:name="variance-pour-uniquement-2-colonnes"/><w:r><w:t xml:space="preserve">Variance (pour uniquement 2 colonnes)</w:t></w:r><w:bookmarkEnd w:id="27"/></w:p><w:p><w:pPr><w:pStyle w:val="SourceCode"/></w:pPr><w:r><w:rPr><w:rStyle w:val="KeywordTok"/></w:rPr><w:t xml:space="preserve">var</w:t></w:r><w:r><w:rPr><w:rStyle w:val="NormalTok"/></w:rPr><w:t xml:space="preserve">(x_matrix[,</w:t></w:r><w:r><w:rPr><w:rStyle w:val="DecValTok"/></w:rPr><w:t xml:space="preserve">1</w:t></w:r><w:r><w:rPr><w:rStyle w:val="OperatorTok"/></w:rPr><w:t xml:space="preserve">:</w:t></w:r><w:r><w:rPr><w:rStyle w:val="DecValTok"/></w:rPr><w:t xml:space="preserve">2</w:t></w:r><w:r><w:rPr><w:rStyle w:val="NormalTok"/></w:rPr><w:t xml:space="preserve">]);</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="SourceCode"/></w:pPr><w:r><w:rPr><w:rStyle w:val="VerbatimChar"/></w:rPr><w:t xml:space="preserve">##                 X971.Guadeloupe X972.Martinique</w:t></w:r><w:r><w:br/></w:r><w:r><w:rPr><w:rStyle w:val="VerbatimChar"/></w:rPr><w:t xml:space="preserve">## X971.Guadeloupe       185185011       147139734</w:t></w:r><w:r><w:br/></w:r><w:r><w:rPr><w:rStyle w:val="VerbatimChar"/></w:rPr><w:t xml:space="preserve">## X972.Martinique       147139734       127959409</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Heading2"/></w:pPr><w:bookmarkStart w:id="28" w:name="correlation-pour-uniquement-2-colonnes"/><w:r><w:t xml:space="preserve">Correlation (pour uniquement 2 colonnes)</w:t></w:r><w:bookmarkEnd w:id="28"/></w:p><w:p><w:pPr><w:pStyle w:val="SourceCode"/></w:pPr><w:r><w:rPr><w:rStyle w:val="KeywordTok"/></w:rPr><w:t xml:space="preserve">cor</w:t></w:r><w:r><w:rPr><w:rStyle w:val="NormalTok"/></w:rPr><w:t xml:space="preserve">(x_matrix[,</w:t></w:r><w:r><w:rPr><w:rStyle w:val="DecValTok"/></w:rPr><w:t xml:space="preserve">1</w:t></w:r><w:r><w:rPr><w:rStyle w:val="OperatorTok"/></w:rPr><w:t xml:space="preserve">:</w:t></w:r><w:r><w:rPr><w:rStyle w:val="DecValTok"/></w:rPr><w:t xml:space="preserve">2</w:t></w:r><w:r><w:rPr><w:rStyle w:val="NormalTok"/></w:rPr><w:t xml:space="preserve">])</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="SourceCode"/></w:pPr><w:r><w:rPr><w:rStyle w:val="VerbatimChar"/></w:rPr><w:t xml:space="preserve">##                 X971.Guadeloupe X972.Martinique</w:t></w:r><w:r><w:br/></w:r><w:r><w:rPr><w:rStyle w:val="VerbatimChar"/></w:rPr><w:t xml:space="preserve">## X971.Guadeloupe       1.0000000       0.9558526</w:t></w:r><w:r><w:br/></w:r><w:r><w:rPr><w:rStyle w:val="VerbatimChar"/></w:rPr><w:t xml:space="preserve">## X972.Martinique       0.9558526       1.0000000</w:t></w:r></w:p>'

$insertionPoint.InsertXML($newSectionsXml)
Write-Output "inserted variance/correlation sections"
